$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(" Oct 30 2020", " Abu Dhabi", "Royals won by 7 wickets (with 15 balls remaining)", "Kings XI Punjab", "Rajasthan Royals", "Chris Gayle ", "99", "63", "6", "8", "157.14"),
    @(" Oct 15 2020", " Sharjah", "Kings XI won by 8 wickets", "Kings XI Punjab", "Royal Challengers Bangalore", "Chris Gayle ", "53", "45", "1", "5", "117.77"),
    @(" Oct 20 2020", " Dubai (DSC)", "Kings XI won by 5 wickets (with 6 balls remaining)", "Kings XI Punjab", "Delhi Capitals", "Chris Gayle ", "29", "13", "3", "2", "223.07"),
    @(" Oct 18 2020", " Dubai (DSC)", "Match tied (Kings XI won the one-over eliminator)", "Kings XI Punjab", "Mumbai Indians", "Chris Gayle ", "24", "21", "1", "2", "114.28"),
    @(" Oct 24 2020", " Dubai (DSC)", "Kings XI won by 12 runs", "Kings XI Punjab", "Sunrisers Hyderabad", "Chris Gayle ", "20", "20", "2", "1", "100.00"),
    @(" Oct 26 2020", " Sharjah", "Kings XI won by 8 wickets (with 7 balls remaining)", "Kings XI Punjab", "Kolkata Knight Riders", "Chris Gayle ", "51", "29", "2", "5", "175.86")
)

# New rows land in A3:K8. Force the whole block to Text format first so
# numeric-looking strings (e.g. "99", "157.14") are written as text instead
# of being auto-converted to numbers by Excel - this matches the source
# file, where every cell (including numeric-looking ones) is t="str".
$targetRange = $ws.Range("A3:K8")
$targetRange.NumberFormat = "@"

$r = 3
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}

# Restore the default "Normal" style so the new cells don't carry an
# explicit style index (matching the source, which has no s="" attributes).
$targetRange.Style = "Normal"
